# ---------------------------------------------------------------------
# Commit: "New crime data collected"
#
# Refreshes the weekly NYPD CompStat precinct report:
#   - Bumps the "Volume 31  Number NN" issue number in the title block.
#   - Rolls the "Report Covering the Week ... Through ..." date range
#     forward by one week.
#   - Overwrites the weekly crime-complaint statistics table (rows
#     15-31: counts + %-change columns for Week-to-Date, 28-Day,
#     Year-to-Date and historical comparisons) with the newly
#     collected figures.
# ---------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a "blank / not-applicable" placeholder into a cell -----
# (displayed as the literal text "0" or "***.*", backed by the
# workbook's existing shared strings) while taking on the number
# format/style of a known-good donor cell that already looks that way.
function Set-TextPlaceholder($ws, $addr, $text, $donorAddr) {
    $dst = $ws.Range($addr)
    $dst.NumberFormat = "@"
    $dst.Value2 = $text
    $donor = $ws.Range($donorAddr)
    $donor.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats
}

# Helper: write a numeric value into a cell that previously held a -----
# text placeholder, restoring a normal numeric style from a donor cell.
function Set-Numeric($ws, $addr, $num, $donorAddr) {
    $dst = $ws.Range($addr)
    $dst.Value2 = $num
    $donor = $ws.Range($donorAddr)
    $donor.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats
}

# --- Title block: "Volume 31   Number  19" -> "...Number  20" ---------
$a8 = $ws.Range("A8")
$a8.Characters(21, 2).Text = "20"

# --- Title block: report week date range -------------------------------
# "Report Covering the Week  5/6/2024  Through  5/12/2024"
#   -> "...5/13/2024  Through  5/19/2024"
# (replace the right-hand token first so the left offset stays valid)
$c9 = $ws.Range("C9")
$c9.Characters(46, 9).Text = "5/19/2024"
$c9.Characters(27, 8).Text = "5/13/2024"

# --- Cells that flip between a numeric value and a text placeholder ---
Set-Numeric $ws "D20" 4 "C15"
Set-Numeric $ws "E20" -75 "L14"
Set-TextPlaceholder $ws "D22" "0" "C14"
Set-TextPlaceholder $ws "E22" "***.*" "E14"
Set-TextPlaceholder $ws "C23" "0" "C14"
Set-Numeric $ws "D23" 1 "C15"
Set-Numeric $ws "E23" -100 "L14"
Set-TextPlaceholder $ws "D28" "0" "C14"
Set-TextPlaceholder $ws "E28" "***.*" "E14"
Set-TextPlaceholder $ws "C31" "0" "C14"
Set-TextPlaceholder $ws "D31" "0" "C14"
Set-TextPlaceholder $ws "E31" "***.*" "E14"

# --- Plain value updates (style/format unchanged) ----------------------
$ws.Range("F15").Value2 = 2
$ws.Range("I15").Value2 = 6
$ws.Range("K15").Value2 = 200
$ws.Range("L15").Value2 = 50
$ws.Range("M15").Value2 = 100
$ws.Range("N15").Value2 = -33.333333333333
$ws.Range("F16").Value2 = 7
$ws.Range("G16").Value2 = 13
$ws.Range("H16").Value2 = -46.153846153846
$ws.Range("I16").Value2 = 40
$ws.Range("J16").Value2 = 55
$ws.Range("K16").Value2 = -27.272727272727
$ws.Range("L16").Value2 = -32.203389830508
$ws.Range("M16").Value2 = 21.212121212121
$ws.Range("N16").Value2 = -82.978723404255
$ws.Range("C17").Value2 = 2
$ws.Range("D17").Value2 = 3
$ws.Range("E17").Value2 = -33.333333333333
$ws.Range("F17").Value2 = 11
$ws.Range("G17").Value2 = 13
$ws.Range("H17").Value2 = -15.384615384615
$ws.Range("I17").Value2 = 44
$ws.Range("J17").Value2 = 54
$ws.Range("K17").Value2 = -18.518518518518
$ws.Range("M17").Value2 = -13.725490196078
$ws.Range("N17").Value2 = -50.561797752809
$ws.Range("F18").Value2 = 4
$ws.Range("G18").Value2 = 10
$ws.Range("H18").Value2 = -60
$ws.Range("I18").Value2 = 34
$ws.Range("J18").Value2 = 44
$ws.Range("K18").Value2 = -22.727272727272
$ws.Range("L18").Value2 = -46.031746031746
$ws.Range("M18").Value2 = -12.820512820512
$ws.Range("N18").Value2 = -82.474226804123
$ws.Range("C19").Value2 = 16
$ws.Range("D19").Value2 = 13
$ws.Range("E19").Value2 = 23.076923076923
$ws.Range("F19").Value2 = 52
$ws.Range("G19").Value2 = 60
$ws.Range("H19").Value2 = -13.333333333333
$ws.Range("I19").Value2 = 252
$ws.Range("J19").Value2 = 282
$ws.Range("K19").Value2 = -10.638297872340
$ws.Range("L19").Value2 = 4.564315352697
$ws.Range("M19").Value2 = 15.068493150684
$ws.Range("N19").Value2 = -21.739130434782
$ws.Range("F20").Value2 = 3
$ws.Range("G20").Value2 = 6
$ws.Range("H20").Value2 = -50
$ws.Range("I20").Value2 = 20
$ws.Range("J20").Value2 = 22
$ws.Range("K20").Value2 = -9.090909090909
$ws.Range("L20").Value2 = -9.090909090909
$ws.Range("M20").Value2 = 17.647058823529
$ws.Range("N20").Value2 = -90.243902439024
$ws.Range("C21").Value2 = 23
$ws.Range("D21").Value2 = 25
$ws.Range("E21").Value2 = -8
$ws.Range("F21").Value2 = 79
$ws.Range("G21").Value2 = 102
$ws.Range("H21").Value2 = -22.549019607843
$ws.Range("I21").Value2 = 396
$ws.Range("J21").Value2 = 459
$ws.Range("K21").Value2 = -13.725490196078
$ws.Range("L21").Value2 = -8.965517241379
$ws.Range("M21").Value2 = 9.090909090909
$ws.Range("N21").Value2 = -62.5
$ws.Range("G22").Value2 = 3
$ws.Range("H22").Value2 = -33.333333333333
$ws.Range("M22").Value2 = 225
$ws.Range("F23").Value2 = 1
$ws.Range("H23").Value2 = -50
$ws.Range("J23").Value2 = 17
$ws.Range("K23").Value2 = -17.647058823529
$ws.Range("L23").Value2 = -44
$ws.Range("M23").Value2 = -26.315789473684
$ws.Range("C24").Value2 = 11
$ws.Range("E24").Value2 = -26.666666666666
$ws.Range("F24").Value2 = 53
$ws.Range("G24").Value2 = 59
$ws.Range("H24").Value2 = -10.169491525423
$ws.Range("I24").Value2 = 233
$ws.Range("J24").Value2 = 258
$ws.Range("K24").Value2 = -9.689922480620
$ws.Range("L24").Value2 = -2.510460251046
$ws.Range("M24").Value2 = -23.856209150326
$ws.Range("C25").Value2 = 10
$ws.Range("D25").Value2 = 12
$ws.Range("E25").Value2 = -16.666666666666
$ws.Range("F25").Value2 = 24
$ws.Range("G25").Value2 = 35
$ws.Range("H25").Value2 = -31.428571428571
$ws.Range("I25").Value2 = 104
$ws.Range("J25").Value2 = 153
$ws.Range("K25").Value2 = -32.026143790849
$ws.Range("L25").Value2 = -12.605042016806
$ws.Range("C26").Value2 = 2
$ws.Range("E26").Value2 = -66.666666666666
$ws.Range("F26").Value2 = 14
$ws.Range("H26").Value2 = -36.363636363636
$ws.Range("I26").Value2 = 115
$ws.Range("J26").Value2 = 127
$ws.Range("K26").Value2 = -9.448818897637
$ws.Range("L26").Value2 = -8
$ws.Range("M26").Value2 = -14.814814814814
$ws.Range("F27").Value2 = 4
$ws.Range("I27").Value2 = 13
$ws.Range("K27").Value2 = 333.333333333333
$ws.Range("L27").Value2 = 160
$ws.Range("C28").Value2 = 3
$ws.Range("F28").Value2 = 4
$ws.Range("G28").Value2 = 6
$ws.Range("H28").Value2 = -33.333333333333
$ws.Range("I28").Value2 = 23
$ws.Range("K28").Value2 = -14.814814814814
$ws.Range("L28").Value2 = 9.523809523809
$ws.Range("F31").Value2 = 2
$ws.Range("G31").Value2 = 1
$ws.Range("H31").Value2 = 100
$ws.Range("I31").Value2 = 3
$ws.Range("K31").Value2 = -72.727272727272
$ws.Range("L31").Value2 = -57.142857142857
